$wb = $excel.ActiveWorkbook

# --- 1. Rename existing headers ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "PO Forecast"

# Copy header formatting (bold, centered, bordered) from an existing header row
$wsWeekly.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# Copy date-cell formatting (custom date/time number format) down column A
$wsWeekly.Range("A2").Copy()
$ws3.Range("A2:A28").PasteSpecial(-4122)

# --- 3. Header row values ---
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# --- 4. Forecast data rows ---
$ws3.Cells.Item(2, 1).Value = 45508.99999999999
$ws3.Cells.Item(2, 2).Value = 0
$ws3.Cells.Item(2, 3).Value = -137.159419619915
$ws3.Cells.Item(2, 4).Value = 129.3021997328305
$ws3.Cells.Item(3, 1).Value = 45515.99999999999
$ws3.Cells.Item(3, 2).Value = 1
$ws3.Cells.Item(3, 3).Value = -135.7230185300907
$ws3.Cells.Item(3, 4).Value = 135.0165540332324
$ws3.Cells.Item(4, 1).Value = 45522.99999999999
$ws3.Cells.Item(4, 2).Value = 10
$ws3.Cells.Item(4, 3).Value = -121.1980591242889
$ws3.Cells.Item(4, 4).Value = 130.905749219192
$ws3.Cells.Item(5, 1).Value = 45529.99999999999
$ws3.Cells.Item(5, 2).Value = 18
$ws3.Cells.Item(5, 3).Value = -117.9286266307527
$ws3.Cells.Item(5, 4).Value = 145.637135790775
$ws3.Cells.Item(6, 1).Value = 45536.99999999999
$ws3.Cells.Item(6, 2).Value = 27
$ws3.Cells.Item(6, 3).Value = -107.340561474685
$ws3.Cells.Item(6, 4).Value = 160.0599156152406
$ws3.Cells.Item(7, 1).Value = 45543.99999999999
$ws3.Cells.Item(7, 2).Value = 36
$ws3.Cells.Item(7, 3).Value = -99.73349147243536
$ws3.Cells.Item(7, 4).Value = 161.8065238984867
$ws3.Cells.Item(8, 1).Value = 45550.99999999999
$ws3.Cells.Item(8, 2).Value = 44
$ws3.Cells.Item(8, 3).Value = -85.67239790553364
$ws3.Cells.Item(8, 4).Value = 175.4695994754048
$ws3.Cells.Item(9, 1).Value = 45557.99999999999
$ws3.Cells.Item(9, 2).Value = 53
$ws3.Cells.Item(9, 3).Value = -77.59215620622869
$ws3.Cells.Item(9, 4).Value = 187.5084702772052
$ws3.Cells.Item(10, 1).Value = 45564.99999999999
$ws3.Cells.Item(10, 2).Value = 61
$ws3.Cells.Item(10, 3).Value = -69.54661997674502
$ws3.Cells.Item(10, 4).Value = 190.2060483559547
$ws3.Cells.Item(11, 1).Value = 45571.99999999999
$ws3.Cells.Item(11, 2).Value = 70
$ws3.Cells.Item(11, 3).Value = -65.15808488094628
$ws3.Cells.Item(11, 4).Value = 205.9908959659287
$ws3.Cells.Item(12, 1).Value = 45578.99999999999
$ws3.Cells.Item(12, 2).Value = 78
$ws3.Cells.Item(12, 3).Value = -52.40234270290221
$ws3.Cells.Item(12, 4).Value = 206.6222508031955
$ws3.Cells.Item(13, 1).Value = 45585.99999999999
$ws3.Cells.Item(13, 2).Value = 87
$ws3.Cells.Item(13, 3).Value = -34.52212794166292
$ws3.Cells.Item(13, 4).Value = 226.2725327110725
$ws3.Cells.Item(14, 1).Value = 45592.99999999999
$ws3.Cells.Item(14, 2).Value = 96
$ws3.Cells.Item(14, 3).Value = -35.15194330528549
$ws3.Cells.Item(14, 4).Value = 224.2489648530013
$ws3.Cells.Item(15, 1).Value = 45606.99999999999
$ws3.Cells.Item(15, 2).Value = 113
$ws3.Cells.Item(15, 3).Value = -22.65417597330101
$ws3.Cells.Item(15, 4).Value = 241.9690764836821
$ws3.Cells.Item(16, 1).Value = 45613.99999999999
$ws3.Cells.Item(16, 2).Value = 121
$ws3.Cells.Item(16, 3).Value = -12.52852142433854
$ws3.Cells.Item(16, 4).Value = 250.2108620187812
$ws3.Cells.Item(17, 1).Value = 45620.99999999999
$ws3.Cells.Item(17, 2).Value = 130
$ws3.Cells.Item(17, 3).Value = -7.130900495275028
$ws3.Cells.Item(17, 4).Value = 260.9356800361904
$ws3.Cells.Item(18, 1).Value = 45627.99999999999
$ws3.Cells.Item(18, 2).Value = 138
$ws3.Cells.Item(18, 3).Value = 14.61113461938557
$ws3.Cells.Item(18, 4).Value = 270.1725775004622
$ws3.Cells.Item(19, 1).Value = 45634.99999999999
$ws3.Cells.Item(19, 2).Value = 147
$ws3.Cells.Item(19, 3).Value = 5.191826127496708
$ws3.Cells.Item(19, 4).Value = 278.0638327516199
$ws3.Cells.Item(20, 1).Value = 45641.99999999999
$ws3.Cells.Item(20, 2).Value = 155
$ws3.Cells.Item(20, 3).Value = 34.84485841360309
$ws3.Cells.Item(20, 4).Value = 281.5847256069842
$ws3.Cells.Item(21, 1).Value = 45648.99999999999
$ws3.Cells.Item(21, 2).Value = 164
$ws3.Cells.Item(21, 3).Value = 34.21038480056657
$ws3.Cells.Item(21, 4).Value = 300.0094934775733
$ws3.Cells.Item(22, 1).Value = 45655.99999999999
$ws3.Cells.Item(22, 2).Value = 173
$ws3.Cells.Item(22, 3).Value = 34.67217404778905
$ws3.Cells.Item(22, 4).Value = 300.6207446994505
$ws3.Cells.Item(23, 1).Value = 45662.99999999999
$ws3.Cells.Item(23, 2).Value = 181
$ws3.Cells.Item(23, 3).Value = 54.30516282335471
$ws3.Cells.Item(23, 4).Value = 310.9795535763253
$ws3.Cells.Item(24, 1).Value = 45669.99999999999
$ws3.Cells.Item(24, 2).Value = 190
$ws3.Cells.Item(24, 3).Value = 64.32236078381044
$ws3.Cells.Item(24, 4).Value = 312.6119314634626
$ws3.Cells.Item(25, 1).Value = 45676.99999999999
$ws3.Cells.Item(25, 2).Value = 198
$ws3.Cells.Item(25, 3).Value = 75.36418612055533
$ws3.Cells.Item(25, 4).Value = 330.75767810639
$ws3.Cells.Item(26, 1).Value = 45683.99999999999
$ws3.Cells.Item(26, 2).Value = 207
$ws3.Cells.Item(26, 3).Value = 81.32231515704079
$ws3.Cells.Item(26, 4).Value = 345.7145992891501
$ws3.Cells.Item(27, 1).Value = 45690.99999999999
$ws3.Cells.Item(27, 2).Value = 215
$ws3.Cells.Item(27, 3).Value = 86.50479941253234
$ws3.Cells.Item(27, 4).Value = 345.0019239686632
$ws3.Cells.Item(28, 1).Value = 45697.99999999999
$ws3.Cells.Item(28, 2).Value = 224
$ws3.Cells.Item(28, 3).Value = 90.31327687776923
$ws3.Cells.Item(28, 4).Value = 342.5726892425683

$excel.CutCopyMode = 0

# Restore the originally active sheet/tab selection
$wsWeekly.Activate()
[void]$wsWeekly.Range("A1").Select()
